# B1-and-B2-PowerPoint.pptx — Sun, Jun 14, 2020 11:04:48 AM
#
# 1) The table on slide 5 ("Type of document / Definition / Why it is
#    important") switches from the custom "Table_0" style to the
#    built-in PowerPoint table style {6D161A82-8FA5-4559-9166-464749602705}.
# 2) The deck's theme is switched from "Integral" (Red Violet colours)
#    to the default "Office Theme" colours.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6D161A82-8FA5-4559-9166-464749602705}")
    }
}

# --- 2. Swap the theme colours to the stock "Office Theme" palette -------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0         # dk1      000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72
